# Update countries & provincias Spain
# - Reorders three pairs of country-name rows (Namibia/Suazilandia,
#   Yemen/Tunez, Montserrat/Islas Malvinas)
# - Refreshes the "last updated" timestamp string
# - Refreshes the COVID-19 stat columns (B:H) for the affected country rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the displayed country name for these three row pairs ---
$ws.Range("A115").Value = "Namibia"      # was Suazilandia
$ws.Range("A116").Value = "Suazilandia"   # was Namibia
$ws.Range("A136").Value = "Yemen"         # was Tunez
$ws.Range("A137").Value = "Tunez"          # was Yemen
$ws.Range("A213").Value = "Montserrat"    # was Islas Malvinas
$ws.Range("A214").Value = "Islas Malvinas" # was Montserrat

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 20:31"

# --- Update statistics values (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 5440779
$ws.Range("C4").Value = 25113
$ws.Range("D4").Value = 2851130
$ws.Range("E4").Value = 2418708
$ws.Range("G4").Value = 526
$ws.Range("H4").Value = 170941

# Row 5 (Brasil)
$ws.Range("B5").Value = 3238216
$ws.Range("C5").Value = 8595
$ws.Range("E5").Value = 775785
$ws.Range("G5").Value = 227
$ws.Range("H5").Value = 105791

# Row 6 (India)
$ws.Range("B6").Value = 2525144
$ws.Range("C6").Value = 65531
$ws.Range("D6").Value = 1805246
$ws.Range("E6").Value = 670764
$ws.Range("G6").Value = 990
$ws.Range("H6").Value = 49134

# Row 12 (Chile)
$ws.Range("B12").Value = 382111
$ws.Range("C12").Value = 2077
$ws.Range("D12").Value = 355037
$ws.Range("E12").Value = 16734
$ws.Range("G12").Value = 41
$ws.Range("H12").Value = 10340

# Row 22 (Alemania)
$ws.Range("B22").Value = 222880
$ws.Range("C22").Value = 611
$ws.Range("E22").Value = 12795
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 9285

# Row 23 (Francia)
$ws.Range("B23").Value = 212211
$ws.Range("C23").Value = 2846
$ws.Range("E23").Value = 98333
$ws.Range("G23").Value = 18
$ws.Range("H23").Value = 30406

# Row 30 (Ecuador)
$ws.Range("B30").Value = 99409
$ws.Range("C30").Value = 1066
$ws.Range("D30").Value = 79176
$ws.Range("E30").Value = 14203
$ws.Range("G30").Value = 20
$ws.Range("H30").Value = 6030

# Row 33 (Israel)
$ws.Range("B33").Value = 91080
$ws.Range("C33").Value = 1258
$ws.Range("D33").Value = 66965
$ws.Range("E33").Value = 23450
$ws.Range("G33").Value = 14
$ws.Range("H33").Value = 665

# Row 41 (Kuwait)
$ws.Range("B41").Value = 75185
$ws.Range("C41").Value = 699
$ws.Range("D41").Value = 66740
$ws.Range("E41").Value = 7951
$ws.Range("G41").Value = 5
$ws.Range("H41").Value = 494

# Row 57 (Marruecos)
$ws.Range("B57").Value = 39241
$ws.Range("C57").Value = 1306
$ws.Range("D57").Value = 27644
$ws.Range("E57").Value = 10986
$ws.Range("G57").Value = 27
$ws.Range("H57").Value = 611

# Row 74 (Chequia)
$ws.Range("B74").Value = 19574
$ws.Range("C74").Value = 173
$ws.Range("D74").Value = 13731
$ws.Range("E74").Value = 5450
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 393

# Row 89 (Zambia)
$ws.Range("B89").Value = 9021
$ws.Range("C89").Value = 358
$ws.Range("D89").Value = 7586
$ws.Range("E89").Value = 1179
$ws.Range("G89").Value = 10
$ws.Range("H89").Value = 256

# Row 95 (Tayikistan)
$ws.Range("B95").Value = 7989
$ws.Range("C95").Value = 39
$ws.Range("D95").Value = 6777
$ws.Range("E95").Value = 1149

# Row 104 (Maldivas)
$ws.Range("B104").Value = 5572
$ws.Range("C104").Value = 78
$ws.Range("D104").Value = 3010
$ws.Range("E104").Value = 2540
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 22

# Row 115 (Suazilandia)
$ws.Range("B115").Value = 3726
$ws.Range("C115").Value = 182
$ws.Range("D115").Value = 848
$ws.Range("E115").Value = 2847
$ws.Range("G115").Value = 4
$ws.Range("H115").Value = 31

# Row 116 (Namibia)
$ws.Range("B116").Value = 3599
$ws.Range("D116").Value = 1991
$ws.Range("E116").Value = 1543
$ws.Range("H116").Value = 65

# Row 125 (Mozambique)
$ws.Range("B125").Value = 2708
$ws.Range("C125").Value = 70
$ws.Range("D125").Value = 1075
$ws.Range("E125").Value = 1614

# Row 136 (Tunez)
$ws.Range("B136").Value = 1858
$ws.Range("C136").Value = 11
$ws.Range("D136").Value = 1009
$ws.Range("E136").Value = 321
$ws.Range("H136").Value = 528

# Row 137 (Yemen)
$ws.Range("D137").Value = 1302
$ws.Range("E137").Value = 492
$ws.Range("H137").Value = 53

# Row 141 (Siria)
$ws.Range("B141").Value = 1515
$ws.Range("C141").Value = 83
$ws.Range("D141").Value = 403
$ws.Range("E141").Value = 1054
$ws.Range("G141").Value = 3
$ws.Range("H141").Value = 58

# Row 144 (Jordania)
$ws.Range("B144").Value = 1329
$ws.Range("C144").Value = 9
$ws.Range("D144").Value = 1229
$ws.Range("E144").Value = 89

# Row 175 (Mauricio)
$ws.Range("B175").Value = 345
$ws.Range("C175").Value = 1
$ws.Range("E175").Value = 1

# Row 213 (Islas Malvinas)
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214 (Montserrat)
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
